# Updated cryptos list on Mon Aug 21 21:17:23 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.181.77"
$ws.Range("E2").Value = "  -1.04%  "
$ws.Range("D3").Value = "1.675.03"
$ws.Range("E3").Value = "  -1.64%  "
$ws.Range("E4").Value = "  -0.81%  "
$ws.Range("D5").Value = "'211.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.61%  "
$ws.Range("D6").Value = "'0.5253"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.97%  "
$ws.Range("E7").Value = "  -0.82%  "
$ws.Range("D8").Value = "'0.2650"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.52%  "
$ws.Range("D9").Value = "'0.06290"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.70%  "
$ws.Range("D10").Value = "'21.33"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.35%  "
$ws.Range("D11").Value = "'0.07557"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.83%  "
$ws.Range("D12").Value = "1.671.19"
$ws.Range("E12").Value = "  -1.96%  "
$ws.Range("D13").Value = "'4.448"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "'0.5606"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.05%  "
$ws.Range("D15").Value = "'66.81"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.20%  "
$ws.Range("D16").Value = "'0.000008012"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.59%  "
$ws.Range("D17").Value = "26.241.35"
$ws.Range("E18").Value = "  -0.79%  "
$ws.Range("D19").Value = "'4.812"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.95%  "
$ws.Range("D20").Value = "'187.70"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.31%  "
$ws.Range("D21").Value = "'10.42"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.20%  "
$ws.Range("D22").Value = "'6.177"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.45%  "
$ws.Range("E23").Value = "  -0.77%  "
$ws.Range("D24").Value = "'149.74"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.34%  "
$ws.Range("D25").Value = "'0.1246"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.96%  "
$ws.Range("D26").Value = "'7.552"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.73%  "
$ws.Range("D27").Value = "'16.02"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.23%  "
$ws.Range("D28").Value = "'0.06154"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.41%  "
$ws.Range("D29").Value = "'1.359"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.61%  "
$ws.Range("D30").Value = "'1.286"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.63%  "
$ws.Range("D31").Value = "'3.492"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.62%  "
$ws.Range("D32").Value = "'3.431"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.88%  "
$ws.Range("D33").Value = "'1.630"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.89%  "
$ws.Range("D34").Value = "'0.9996"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.34%  "
$ws.Range("D35").Value = "'0.6059"
$ws.Range("D35").Style = "Normal"
$ws.Range("E36").Value = "  -0.33%  "
$ws.Range("D37").Value = "'2.737"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.26%  "
$ws.Range("D38").Value = "'6.088"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.01%  "
$ws.Range("D39").Value = "'0.01612"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.16%  "
$ws.Range("D40").Value = "1.082.55"
$ws.Range("E40").Value = "  -3.41%  "
$ws.Range("D41").Value = "'0.8697"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.57%  "
$ws.Range("E42").Value = "  -1.16%  "
$ws.Range("D43").Value = "'99.97"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.60%  "
$ws.Range("D44").Value = "1.826.30"
$ws.Range("E44").Value = "  -1.59%  "
$ws.Range("E45").Value = "  +0.80%  "
$ws.Range("D46").Value = "'56.00"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.11%  "
$ws.Range("D47").Value = "'0.9974"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.88%  "
$ws.Range("D48").Value = "'8.023"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.68%  "
$ws.Range("D49").Value = "'0.05234"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.10%  "
$ws.Range("D50").Value = "'0.4257"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.18%  "
$ws.Range("D51").Value = "'5.965"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.92%  "
